$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.6881615054955936
$ws.Range("D2").Value = 0.6264699153511782
$ws.Range("E2").Value = 0.6221844651985571
$ws.Range("G2").Value = 0.7351167062838367
$ws.Range("H2").Value = 0.6865248226950355
$ws.Range("I2").Value = 0.6844945096892918
$ws.Range("K2").Value = 0.7148942439039832
$ws.Range("L2").Value = 0.6844200411805079
$ws.Range("M2").Value = 0.6839613283381253
$ws.Range("O2").Value = 0.7926823995174128
$ws.Range("P2").Value = 0.7789750629146649
$ws.Range("Q2").Value = 0.7790072107894227
$ws.Range("S2").Value = 0.7708159889374424
$ws.Range("T2").Value = 0.7126058110272249
$ws.Range("U2").Value = 0.715969873474698
$ws.Range("W2").Value = 0.7730195373019669
$ws.Range("X2").Value = 0.7147334706016929
$ws.Range("Y2").Value = 0.7213536650873087
$ws.Range("C3").Value = 0.7985956243374288
$ws.Range("D3").Value = 0.7875543353923588
$ws.Range("E3").Value = 0.7873394143488996
$ws.Range("G3").Value = 0.8128679272287336
$ws.Range("H3").Value = 0.8004118050789293
$ws.Range("I3").Value = 0.801663318551846
$ws.Range("K3").Value = 0.8219350976610016
$ws.Range("L3").Value = 0.8133150308853809
$ws.Range("M3").Value = 0.8150161887268537
$ws.Range("O3").Value = 0.8391740644712609
$ws.Range("P3").Value = 0.8305193319606496
$ws.Range("Q3").Value = 0.8300626228518213
$ws.Range("S3").Value = 0.8595560068448369
$ws.Range("T3").Value = 0.8519103180050331
$ws.Range("U3").Value = 0.8510766315601531
$ws.Range("W3").Value = 0.8657432086529147
$ws.Range("X3").Value = 0.8561885152139099
$ws.Range("Y3").Value = 0.8559497519229347
$ws.Range("C4").Value = 0.8231285042630135
$ws.Range("D4").Value = 0.8153283001601463
$ws.Range("E4").Value = 0.8145030727138709
$ws.Range("G4").Value = 0.8611739633378448
$ws.Range("H4").Value = 0.8540150995195608
$ws.Range("I4").Value = 0.854172320614113
$ws.Range("K4").Value = 0.8398044924418351
$ws.Range("L4").Value = 0.8325326012354154
$ws.Range("M4").Value = 0.8322725549148501
$ws.Range("O4").Value = 0.838719606885767
$ws.Range("P4").Value = 0.8347517730496454
$ws.Range("Q4").Value = 0.8340428381639159
$ws.Range("S4").Value = 0.8524934376831117
$ws.Range("T4").Value = 0.8476549988560971
$ws.Range("U4").Value = 0.8469769510436667
$ws.Range("W4").Value = 0.8545729195418804
$ws.Range("X4").Value = 0.8498055364905056
$ws.Range("Y4").Value = 0.8490848129987272
$ws.Range("C5").Value = 0.8147226913677754
$ws.Range("D5").Value = 0.8047357584076871
$ws.Range("E5").Value = 0.8026002501243894
$ws.Range("G5").Value = 0.8478250140868472
$ws.Range("H5").Value = 0.8390528483184626
$ws.Range("I5").Value = 0.8388895481124227
$ws.Range("O5").Value = 0.832594934784729
$ws.Range("P5").Value = 0.8197208876687258
$ws.Range("Q5").Value = 0.8173006130378528
$ws.Range("S5").Value = 0.8567157216617579
$ws.Range("T5").Value = 0.8455044612216884
$ws.Range("U5").Value = 0.8447762295277957
$ws.Range("C6").Value = 0.718150479358624
$ws.Range("D6").Value = 0.6844886753603293
$ws.Range("E6").Value = 0.678974110704153
$ws.Range("G6").Value = 0.7488661302646362
$ws.Range("H6").Value = 0.7059254175245938
$ws.Range("I6").Value = 0.7033312625658155
$ws.Range("K6").Value = 0.7520611880317394
$ws.Range("L6").Value = 0.7037977579501258
$ws.Range("M6").Value = 0.6995047456688194
$ws.Range("O6").Value = 0.7359756673928738
$ws.Range("P6").Value = 0.7145504461221688
$ws.Range("Q6").Value = 0.7078176161253705
$ws.Range("S6").Value = 0.7455727442185192
$ws.Range("T6").Value = 0.7188515213909861
$ws.Range("U6").Value = 0.712972901679728
$ws.Range("W6").Value = 0.7520611880317394
$ws.Range("X6").Value = 0.7037977579501258
$ws.Range("Y6").Value = 0.6995047456688194
